# Apply "added some more cues" edit to the Russman cues sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shout_Grenade_Warning_Lure: track becomes a list (was a lone number),
#     and the old "Monkey bomb" description note is removed ---
$ws.Range("B11").Value = "15, 280"
$ws.Range("D11").ClearContents()

# --- Shout_PickAxePowerAttack: prepend track 39 to the existing list ---
$ws.Range("B8").Value = "39, 241, 242, 243, 244, 245, 246, 247, 248"

# --- Shout_LaserPoint_Generic: add a track number ---
$ws.Range("B18").Value = "55"

# --- Two brand-new cue rows appended after the existing data ---
$ws.Range("B49").Value = "64"
$ws.Range("D49").Value = "ought to be able to do something with this"

$ws.Range("B50").Value = "68"
$ws.Range("D50").Value = "can't keep hauling this around"

# --- Shout_FlyingGrabber: add a track number ---
$ws.Range("B4").Value = "74"

# --- Shout_ProtectDrilldozer_WhenTakingDamage: extend the track list ---
$ws.Range("B26").Value = "89, 90, 97, 98, 99, 100"

# Match the author's final selection in the sheet
$ws.Range("B26").Select()
